$wb = $excel.ActiveWorkbook

# Rows and new values that changed for F column on the "展览" and "全部类型" sheets.
$updates = @{
    2  = 1093
    5  = 4662
    8  = 1396
    11 = 1143
    13 = 633
    15 = 36
    16 = 15
    17 = 275
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
